$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.518.44'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.627.70'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '213.18'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '0.504'
$ws.Range('E6').Value = '  +2.18%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '0.0623'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '18.80'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '1.853.64'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.14'
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.592.92'
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '65.12'
$ws.Range('E16').Value = '  +3.47%  '
$ws.Range('D17').Value = '26.538.25'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '214.80'
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +12.41%  '
$ws.Range('D25').Value = '147.51'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +2.20%  '
$ws.Range('D29').Value = '15.57'
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D33').Value = '2.96'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').Value = '1.239.39'
$ws.Range('E35').Value = '  +6.06%  '
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('E37').Value = '  +4.69%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('D42').Value = '0.799'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').Value = '1.764.06'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '93.15'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('D47').Value = '54.86'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.52'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('E51').Value = '  -0.54%  '
